# Add support for specifying min, max and default value for numeric
# ("open-num") questions by populating the Options column (D) for the
# two existing open-num rows:
#   Row 4 -> "How many years have you been using R?"  (min=1, max=27,  default=3)
#   Row 3 -> "How old are you?"                        (min=1, max=100, default=20)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("questions")

# Write row 4 first so the new shared-string table entries are created in
# the same order as the source edit ("1,27,3" before "1,100,20").
$ws.Range("D4").Value = "1,27,3"
$ws.Range("D3").Value = "1,100,20"

# Leave the selection on the last-edited cell, as in the original edit.
$ws.Range("D4").Select()
